$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The E2 cell held a multi-line poster-author list with embedded line
# breaks; replace it with the same text joined onto a single line (each
# line break collapsed to a single space).
$newText = "Zvonimir Banoža, Matija Radobuljac, Sanja Pavić Jelečki, Suzana Palatinuš, Ivana Kralj, Saša Balija, Silvia Tisaj Pigac, Ivana Sklepić Klobučarić, Tina Kresonja, Milena Škvorc, Suzana Palatinuš, Barbara Samvik, Kristijan Štefanes, Nikola Čopor, Ivica Bračko, Ramon Tumbas, Emil Kralj, Dragica Svetličić"

$ws.Range("E2").Value = $newText

# Move/collapse the active selection to the single cell E2 (was F2 with
# the whole row 2 selected).
$ws.Range("E2").Select()
